# Rename "saitama_temp_202208" to "saitama_temp_202308" and refresh its
# data with the August-2023 readings (the sheet previously held August 2022
# data). Also move the sheet's selection to the full data block A2:D32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("saitama_temp_202208")

$ws.Name = "saitama_temp_202308"

# New values for columns A (date serial), B (avg temp), C (max temp),
# D (min temp) for rows 2..32 (2023-08-01 .. 2023-08-31).
$newData = @(
    @(45139, 26.9, 31, 23.8),
    @(45140, 28.9, 35.9, 23),
    @(45141, 30.2, 35.9, 25),
    @(45142, 31.5, 37.9, 25.7),
    @(45143, 31, 36.6, 26.9),
    @(45144, 29.7, 34.1, 27),
    @(45145, 30.1, 35, 26.6),
    @(45146, 29.2, 35.200000000000003, 25.4),
    @(45147, 28.3, 32.9, 26),
    @(45148, 29.8, 34.6, 26.5),
    @(45149, 28.9, 34, 25),
    @(45150, 30, 35.6, 25.5),
    @(45151, 28, 32.1, 25.3),
    @(45152, 28.7, 32.299999999999997, 26.4),
    @(45153, 29, 32.5, 26.1),
    @(45154, 28.7, 32, 26.1),
    @(45155, 30.1, 35.299999999999997, 25.9),
    @(45156, 30.2, 35.799999999999997, 26.6),
    @(45157, 30.2, 36, 26.3),
    @(45158, 30.7, 35.6, 26.1),
    @(45159, 30.1, 34.299999999999997, 26),
    @(45160, 28.2, 32.200000000000003, 26.3),
    @(45161, 29.4, 33.9, 26.7),
    @(45162, 29.8, 34, 26.1),
    @(45163, 30.1, 35.299999999999997, 25.9),
    @(45164, 29, 33.299999999999997, 25.8),
    @(45165, 27.5, 31.9, 24.1),
    @(45166, 29.8, 35.5, 25.2),
    @(45167, 29.8, 36.299999999999997, 25.1),
    @(45168, 29.8, 35.1, 25.7),
    @(45169, 30, 35, 26.9)
)

$r = 2
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$ws.Activate()
$ws.Range("A2:D32").Select()
